$wb = $excel.ActiveWorkbook

# Rename the first sheet from "HUP99-complete" to "HUP099-complete"
# (zero-padding the patient id to 3 digits).
$ws = $wb.Worksheets.Item("HUP99-complete")
$ws.Name = "HUP099-complete"
